$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Advance the date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = $ws.Range("A1").Value2 + 1

# Update prices
$ws.Range("D29").Value = 410
$ws.Range("D30").Value = 445
